# Update the "Förändrad" (Changed) date column (C) for rows 2-10
# from 2023-10-25 (45224) to 2023-11-03 (45233).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 10; $row++) {
    $cell = $ws.Range("C$row")
    if ($cell.Value2 -eq 45224) {
        $cell.Value = 45233
    }
}
